$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BOSS_HPBAR_BACK block (rows 7-12): Position_X/Y, Size_X sync fix
$ws.Range("B9").Value = 977
$ws.Range("B10").Value = 35
$ws.Range("B11").Value = 427

# PLAYER_HPBAR_FRONT block (rows 31-36): Position_X, Size_X sync fix
$ws.Range("B33").Value = 139
$ws.Range("B35").Value = 283

# Restore default view: scroll back to top-left, select E13
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
